$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 211.3819611842011
$ws.Range("C2").Value = 51.65877335318159

$ws.Range("B7").Value = 79268.23544407543
$ws.Range("C7").Value = 19372.04000744309

$ws.Range("B8").Value = 1189.023531661131
$ws.Range("C8").Value = 290.5806001116464

$ws.Range("B9").Value = 52845.49029605028
$ws.Range("C9").Value = 12914.6933382954
